# Regenerate the "K" column (column G) of save_data using the strike
# information (previously "Strike#") instead of the old values, after
# recomputing std/mean and writing the resulting s_vals.
#
# The workbook only has one worksheet holding the table; column G (header
# "K") is updated row by row to the freshly computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G ("K"), as produced by the
# regenerated std/mean calculation (s_vals).
$kValues = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 1
    6  = 0
    7  = 0
    8  = 1
    9  = 1
    10 = 2
    11 = 0
    12 = 0
    13 = 1
    14 = 2
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 1
    23 = 1
    24 = 1
    25 = 2
    26 = 1
    27 = 1
    28 = 1
    29 = 0
    30 = 1
    31 = 0
    32 = 1
    33 = 0
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 0
    39 = 1
    40 = 1
    41 = 1
    42 = 1
    43 = 1
    44 = 1
    45 = 1
    46 = 0
    47 = 2
    48 = 0
    51 = 1
    53 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
